# Refresh the cryptos price/volume snapshot (GitHub Actions bot update).
#
# Most cells are non-numeric-looking strings ("1.641.82", "  +1.20%  ", …)
# and a plain .Value assignment stores them as text, same as the source
# file. A handful of Price cells (D5, D11, D16, D19, D24, D25, D26, D29,
# D45, D46, D50) look like ordinary decimals (e.g. "217.36"), so Excel
# would otherwise auto-convert them to numbers; those are written with a
# leading apostrophe to keep them text, then snapped back to the Normal
# style so no stray number-format/quote-prefix styling is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.198.63"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "1.642.75"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'217.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("D11").Value = "'0.0849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "1.871.98"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "1.649.52"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("D16").Value = "'67.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "27.176.38"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("D19").Value = "'218.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  +3.82%  "
$ws.Range("E22").Value = "  +7.12%  "
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "'9.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "'147.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").Value = "'7.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").Value = "'15.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").Value = "1.275.07"
$ws.Range("E35").Value = "  +2.44%  "
$ws.Range("E37").Value = "  +1.86%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "1.782.29"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "'61.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("D46").Value = "'91.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").Value = "'7.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("E51").Value = "  +0.11%  "
